$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.21%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.25%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.918"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.59%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07320"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.249"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'22.73%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.45%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-1.25%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9043"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.83%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'18.65%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1687"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.72%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08270"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.47%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03123"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.99%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09926"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.88%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.00%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005716"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.91%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.523"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.61%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.047"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.59%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3329"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.37%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.72%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.201"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.74%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-12.04%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04513"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.42%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004156"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-10.80%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.95%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003398"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-95.46%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01571"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.68%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04437"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.88%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.79%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009565"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-4.97%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.65%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'13.70%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008339"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-9.29%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006110"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.18%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.277"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-1.44%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.02%"
$ws.Range("E51").Style = "Normal"
